$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("I8").Value = "test"
$ws.Range("I8").Font.Underline = $true
$ws.Range("I8").Font.Color = 15570276
